# Update database and shift "read_price" yearly columns by one year,
# adding the new 1401/12 period and dropping the oldest 1396/12 period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows: shift the twelve-month period labels left by one year ---
# (E<-old F, F<-old G, G<-old H, H<-old I, I<-new 1401/12 period)
$headerRows = @(8, 24)
foreach ($r in $headerRows) {
    $ws.Range("E$r").Value = "دوازده ماهه منتهی به 1397/12"
    $ws.Range("F$r").Value = "دوازده ماهه منتهی به 1398/12"
    $ws.Range("G$r").Value = "دوازده ماهه منتهی به 1399/12"
    $ws.Range("H$r").Value = "دوازده ماهه منتهی به 1400/12"
    $ws.Range("I$r").Value = "دوازده ماهه منتهی به 1401/12"
}

# --- Data table: هزینه های عمومی و اداری (general & administrative expenses) ---
# Each row's values shift one column left (E<-F, F<-G, G<-H, H<-I) and the
# new I column is populated with the freshly reported 1401/12 figure.

# هزینه حمل و نقل و انتقال
$ws.Range("E10").Value = 1514801
$ws.Range("F10").Value = 2488502
$ws.Range("G10").Value = 8042261
$ws.Range("H10").Value = 7702800
$ws.Range("I10").Value = 5756768

# هزینه خدمات پس از فروش
$ws.Range("E11").Value = 19707
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0

# حق العمل و کمیسیون فروش (unchanged, still all zero)
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

# هزینه تبلیغات (unchanged, still all zero)
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# هزینه مواد مصرفی
$ws.Range("E14").Value = 9754
$ws.Range("F14").Value = 9168
$ws.Range("G14").Value = 23259
$ws.Range("H14").Value = 104540
$ws.Range("I14").Value = 171999

# هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Range("E15").Value = 1513
$ws.Range("F15").Value = 2092
$ws.Range("G15").Value = 1697
$ws.Range("H15").Value = 2944
$ws.Range("I15").Value = 2317

# هزینه استهلاک
$ws.Range("E16").Value = 7500
$ws.Range("F16").Value = 7500
$ws.Range("G16").Value = 7641
$ws.Range("H16").Value = 28612
$ws.Range("I16").Value = 46832

# هزینه حقوق و دستمزد
$ws.Range("E17").Value = 98776
$ws.Range("F17").Value = 123594
$ws.Range("G17").Value = 328061
$ws.Range("H17").Value = 429496
$ws.Range("I17").Value = 1058985

# هزینه مطالبات مشکوک الوصول
$ws.Range("E18").Value = 677609
$ws.Range("F18").Value = 298881
$ws.Range("G18").Value = 277514
$ws.Range("H18").Value = 6191
$ws.Range("I18").Value = 505418

# سایر هزینه ها
$ws.Range("E19").Value = 818695
$ws.Range("F19").Value = 653717
$ws.Range("G19").Value = 931446
$ws.Range("H19").Value = 1482762
$ws.Range("I19").Value = 3171277

# جمع (total)
$ws.Range("E20").Value = 3148355
$ws.Range("F20").Value = 3583454
$ws.Range("G20").Value = 9611879
$ws.Range("H20").Value = 9757345
$ws.Range("I20").Value = 10713596

# --- تعداد پرسنل (personnel counts) ---

# تعداد پرسنل غیر تولیدی شرکت
$ws.Range("E26").Value = 961
$ws.Range("F26").Value = 966
$ws.Range("G26").Value = 628
$ws.Range("H26").Value = 648
$ws.Range("I26").Value = 542

# تعداد پرسنل تولیدی شرکت
$ws.Range("E27").Value = 186
$ws.Range("F27").Value = 188
$ws.Range("G27").Value = 532
$ws.Range("H27").Value = 695
$ws.Range("I27").Value = 877

$wb.Save()
